$d = $word.ActiveDocument

# --- Part 1: fix "parallelisable" -> "parallélisable" and relocate the
#     "_GoBack" bookmark so it now sits between the new "é" run and the
#     trailing "lisable" run.

$r = $d.Content
$r.Find.Execute("parallelisable") | Out-Null
$start = $r.Start

# Replace the bare "e" (offset 6 within the word) with an accented "é".
$eRange = $d.Range($start + 6, $start + 7)
$eRange.Text = "é"

# Force a run boundary between "parall" and "é" using a scratch bookmark,
# then drop the scratch bookmark - the run split survives its removal.
$splitRange = $d.Range($start + 6, $start + 6)
$d.Bookmarks.Add("TempSplit", $splitRange) | Out-Null

# Re-create "_GoBack" right after the "é", before "lisable".
$goBackRange = $d.Range($start + 7, $start + 7)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null

$d.Bookmarks("TempSplit").Delete()

# --- Part 2: merge the old "Couleur" / " dessin : " runs (which used to
#     carry the "_GoBack" bookmark) back into a single run. We capture the
#     live text first (rather than retyping it) so characters such as the
#     non-breaking space before the colon are preserved exactly.

$r2 = $d.Content
$r2.Find.MatchCase = $true
$r2.Find.Execute("Couleur") | Out-Null
$anchorStart = $r2.Start

$mergedRange = $d.Range($anchorStart, $anchorStart + 17)
$mergedText = $mergedRange.Text

$d.Content.Find.Execute($mergedText, $true, $false, $false, $false, `
    $false, $true, 1, $false, $mergedText, 2) | Out-Null
